# Generate Report for Archive
#
# - Update the localization status text "Ready for handoff" -> "In Translation"
#   everywhere it appears (Overview!E2/F2, zh-cn!C2, de-de!C2 - all backed by
#   the same shared string).
# - Narrow the status column(s) to match the new, shorter text's autofit width
#   (Overview columns E & F, and the "Status" column (C) on the zh-cn / de-de
#   detail sheets).

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# ---- zh-cn detail sheet -----------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

# ---- de-de detail sheet -----------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
